$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.422.21'
$ws.Range('E2').Value = '  -5.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.224.06'
$ws.Range('E3').Value = '  -5.89%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.91'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.10'
$ws.Range('E7').Value = '  -5.69%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.560'
$ws.Range('E9').Value = '  -6.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.19'
$ws.Range('E10').Value = '  +5.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0958'
$ws.Range('E11').Value = '  -6.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.34'
$ws.Range('E12').Value = '  -2.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('E14').Value = '  -6.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.553.80'
$ws.Range('E15').Value = '  -6.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.87'
$ws.Range('E16').Value = '  -9.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.844'
$ws.Range('E17').Value = '  -8.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.224.53'
$ws.Range('E18').Value = '  -6.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.315.86'
$ws.Range('E19').Value = '  -5.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0958'
$ws.Range('E20').Value = '  -7.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.40'
$ws.Range('E21').Value = '  -6.26%  '
$ws.Range('E22').Value = '  -7.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.56'
$ws.Range('E23').Value = '  -8.36%  '
$ws.Range('E24').Value = '  +13.00%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.63'
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.87'
$ws.Range('E28').Value = '  -7.11%  '
$ws.Range('E29').Value = '  -4.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.43'
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.60'
$ws.Range('E31').Value = '  -7.89%  '
$ws.Range('E32').Value = '  -7.42%  '
$ws.Range('E33').Value = '  -6.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0719'
$ws.Range('E34').Value = '  -4.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.26'
$ws.Range('E35').Value = '  -3.71%  '
$ws.Range('E36').Value = '  -9.62%  '
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.25'
$ws.Range('E38').Value = '  +16.72%  '
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.27'
$ws.Range('E40').Value = '  -4.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.87'
$ws.Range('E41').Value = '  -11.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.95'
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('E43').Value = '  +3.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.99'
$ws.Range('E44').Value = '  -12.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.97'
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.101'
$ws.Range('E46').Value = '  -6.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.83'
$ws.Range('E47').Value = '  +11.06%  '
$ws.Range('E48').Value = '  +6.47%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').Value = '  -5.52%  '
$ws.Range('E51').Value = '  -5.07%  '
